$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Text fix-ups in the existing "第十二周周三" (week 12, Wed) blocks:
#    "完成聊天" -> "完成单聊聊天"
# ------------------------------------------------------------------
$ws.Range("B200").Value = "完成单聊聊天"
$ws.Range("B210").Value = "完成单聊聊天"

# ------------------------------------------------------------------
# 2. Newly filled-in completion percentages in the second block
#    (rows 204-212, "日期：2018.11.26 第十三周周一")
# ------------------------------------------------------------------
$ws.Range("C207").Value = 0.7
$ws.Range("C208").Value = 1
$ws.Range("C210").Value = 1

# ------------------------------------------------------------------
# 3. Append a brand-new weekly block in rows 214-222, mirroring the
#    layout/formatting of the rows 204-212 block immediately above it.
#    Copy *formats only* so no new cell styles are introduced.
# ------------------------------------------------------------------
$src = $ws.Range("A204:D212")
$dst = $ws.Range("A214:D222")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Title row (merged A214:D214)
$ws.Range("A214").Value = "日期：2018.11.28 第十三周周三"

# Header row
$ws.Range("A215").Value = "组员"
$ws.Range("B215").Value = "计划内容"
$ws.Range("C215").Value = "完成情况"
$ws.Range("D215").Value = "备注"

# Member rows
$ws.Range("A216").Value = "王伟锋"
$ws.Range("B216").Value = "基本功能已完成，协助app完成"
$ws.Range("D216").Value = "协助情况不计入任务完成情况"

$ws.Range("A217").Value = "陈升云"
$ws.Range("B217").Value = "完成消息页面的优化和各功能的实现"

$ws.Range("A218").Value = "林玮成"
$ws.Range("B218").Value = "等待最终程序的测试"

$ws.Range("A219").Value = "吴帅辰"
$ws.Range("B219").Value = "基本功能已完成，协助app完成"
$ws.Range("D219").Value = "协助情况不计入任务完成情况"

$ws.Range("A220").Value = "李海洋"
$ws.Range("B220").Value = "完成群聊"

# Summary row (merged A221:D221 originally, now merged A221:D222 with the
# trailing blank row)
$ws.Range("A221").Value = "总结："

# Merge the new title and summary rows like every other block on the sheet
$ws.Range("A214:D214").Merge()
$ws.Range("A221:D222").Merge()

# ------------------------------------------------------------------
# 4. Move the active selection to match the post-edit state
# ------------------------------------------------------------------
$ws.Range("F216").Select()
